$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2016365
$ws.Range("I38").Value = 2304224.2
$ws.Range("J38").Value = 1350
$ws.Range("K38").Value = 6912672.600000001
$ws.Range("L38").Value = 4050
$ws.Range("M38").Value = -6912300.600000001
$ws.Range("N38").Value = -4794

$ws.Range("H76").Value = 4340.0625
$ws.Range("I76").Value = 4226.077
$ws.Range("K76").Value = 4226.077
$ws.Range("M76").Value = -3911.077

$ws.Range("H79").Value = 4340.0625
$ws.Range("I79").Value = 4226.077
$ws.Range("K79").Value = 4226.077
$ws.Range("M79").Value = -3134.077

$ws.Range("H137").Value = 1156.7174
$ws.Range("I137").Value = 848.8919
$ws.Range("J137").Value = 2422.2222
$ws.Range("K137").Value = 2546.6757
$ws.Range("L137").Value = 7266.6666
$ws.Range("M137").Value = 3.324300000000221
$ws.Range("N137").Value = -12366.6666

$ws.Range("H138").Value = 2199.32
$ws.Range("I138").Value = 1544.091
$ws.Range("J138").Value = 2280.3035
$ws.Range("K138").Value = 4632.272999999999
$ws.Range("L138").Value = 6840.9105
$ws.Range("M138").Value = 507.7270000000008
$ws.Range("N138").Value = -17120.9105

$ws.Range("H141").Value = 2576.2593
$ws.Range("I141").Value = 2177.842
$ws.Range("K141").Value = 6533.526
$ws.Range("M141").Value = -1353.526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3368.98
$ws.Range("I32").Value = 2748.6924
$ws.Range("J32").Value = 5568.1816
$ws.Range("K32").Value = 2748.6924
$ws.Range("L32").Value = 5568.1816
$ws.Range("M32").Value = -2461.6924
$ws.Range("N32").Value = -6142.1816

$ws.Range("H110").Value = 83508800
$ws.Range("I110").Value = 91100450
$ws.Range("J110").Value = 700
$ws.Range("K110").Value = 91100450
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = -91098405
$ws.Range("N110").Value = -4790

$ws.Range("H132").Value = 12172.107
$ws.Range("I132").Value = 15016.768
$ws.Range("J132").Value = 2762.8462
$ws.Range("K132").Value = 45050.304
$ws.Range("L132").Value = 8288.5386
$ws.Range("M132").Value = -42520.304
$ws.Range("N132").Value = -13348.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 139925.75
$ws.Range("I86").Value = 185901
$ws.Range("K86").Value = 185901
$ws.Range("M86").Value = -184778

$ws.Range("H89").Value = 139925.75
$ws.Range("I89").Value = 185901
$ws.Range("K89").Value = 929505
$ws.Range("M89").Value = -923889

$ws.Range("H107").Value = 66698196
$ws.Range("I107").Value = 71462250
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 71462250
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -71460330
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 5247.5
$ws.Range("I12").Value = 990
$ws.Range("J12").Value = 6666.6665
$ws.Range("K12").Value = 990
$ws.Range("L12").Value = 6666.6665
$ws.Range("M12").Value = -820
$ws.Range("N12").Value = -7006.6665

$ws.Range("H31").Value = 40673.758
$ws.Range("I31").Value = 980.5789
$ws.Range("J31").Value = 82572.11
$ws.Range("K31").Value = 980.5789
$ws.Range("L31").Value = 82572.11
$ws.Range("M31").Value = -685.5789
$ws.Range("N31").Value = -83162.11

$ws.Range("H34").Value = 40673.758
$ws.Range("I34").Value = 980.5789
$ws.Range("J34").Value = 82572.11
$ws.Range("K34").Value = 980.5789
$ws.Range("L34").Value = 82572.11
$ws.Range("M34").Value = -778.5789
$ws.Range("N34").Value = -82976.11

$ws.Range("H94").Value = 1590.1111
$ws.Range("I94").Value = 1606
$ws.Range("J94").Value = 1585.5714
$ws.Range("K94").Value = 1606
$ws.Range("L94").Value = 1585.5714
$ws.Range("M94").Value = -1155
$ws.Range("N94").Value = -2487.5714

$ws.Range("H132").Value = 3182.3044
$ws.Range("I132").Value = 3326.0557
$ws.Range("J132").Value = 2664.8
$ws.Range("K132").Value = 9978.167099999999
$ws.Range("L132").Value = 7994.400000000001
$ws.Range("M132").Value = -7448.167099999999
$ws.Range("N132").Value = -13054.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 62.736843
$ws.Range("I15").Value = 30.625
$ws.Range("J15").Value = 234
$ws.Range("K15").Value = 91.875
$ws.Range("L15").Value = 702
$ws.Range("M15").Value = 48.125
$ws.Range("N15").Value = -982

$ws.Range("H75").Value = 1463.3334
$ws.Range("J75").Value = 2243.3333
$ws.Range("L75").Value = 6729.999899999999
$ws.Range("N75").Value = -8725.999899999999

$ws.Range("H78").Value = 1463.3334
$ws.Range("J78").Value = 2243.3333
$ws.Range("L78").Value = 20189.9997
$ws.Range("N78").Value = -30173.9997

$ws.Range("H107").Value = 1131.0714
$ws.Range("I107").Value = 786.6667
$ws.Range("J107").Value = 1225
$ws.Range("K107").Value = 2360.0001
$ws.Range("L107").Value = 3675
$ws.Range("M107").Value = -440.0001000000002
$ws.Range("N107").Value = -7515

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3679473.8
$ws.Range("I126").Value = 3208.3635
$ws.Range("J126").Value = 11767257
$ws.Range("K126").Value = 9625.0905
$ws.Range("L126").Value = 35301771
$ws.Range("M126").Value = -7155.0905
$ws.Range("N126").Value = -35306711

$ws.Range("H132").Value = 2603.55
$ws.Range("I132").Value = 1841.5938
$ws.Range("J132").Value = 5651.375
$ws.Range("K132").Value = 5524.7814
$ws.Range("L132").Value = 16954.125
$ws.Range("M132").Value = -2994.7814
$ws.Range("N132").Value = -22014.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7550
$ws.Range("I46").Value = 401
$ws.Range("K46").Value = 401
$ws.Range("M46").Value = -213

$ws.Range("H125").Value = 47995
$ws.Range("J125").Value = 47995
$ws.Range("L125").Value = 47995
$ws.Range("N125").Value = -57835

$ws.Range("H132").Value = 3268.2903
$ws.Range("I132").Value = 3353.2173
$ws.Range("J132").Value = 3024.125
$ws.Range("K132").Value = 10059.6519
$ws.Range("L132").Value = 9072.375
$ws.Range("M132").Value = -7529.651899999999
$ws.Range("N132").Value = -14132.375

$ws.Range("H136").Value = 2164.4856
$ws.Range("I136").Value = 1537.8334
$ws.Range("K136").Value = 4613.5002
$ws.Range("M136").Value = -2063.5002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4417.6807
$ws.Range("I132").Value = 2305.4243
$ws.Range("J132").Value = 9396.571
$ws.Range("K132").Value = 6916.2729
$ws.Range("L132").Value = 28189.713
$ws.Range("M132").Value = -4386.2729
$ws.Range("N132").Value = -33249.713
